$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text and the company name value (typo fix from original data)
$ws.Range("D1").Value = "nomeDaEmpresa"
$ws.Range("D2").Value = "world app company"

# Remove the stray underline formatting that had been applied to A1 and D2
$ws.Range("A1").Font.Underline = -4142
$ws.Range("D2").Font.Underline = -4142

# Move the active selection to E2
$ws.Range("E2").Select()
